# Run Manager.xlsx - "Featured added to invoke methods Dynamically"
#
# Inserts a new "Test Scenario" column at the front of the "Test Info" sheet
# (with value "LoginPageTest" for the existing test row), and updates the
# active-sheet / selection view state so that "Test Info" is the selected tab
# and "Test Data" keeps a plain (non-tab-selected) view with its selection
# moved to F3.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Test Info"
$ws2 = $wb.Worksheets.Item(2)   # "Test Data"

# --- Shift existing columns A:E of "Test Info" one column to the right (B:F),
#     preserving their values, then put the new "Test Scenario" data in A. ---
$headerRow = @(
    $ws1.Range("A1").Value2,
    $ws1.Range("B1").Value2,
    $ws1.Range("C1").Value2,
    $ws1.Range("D1").Value2,
    $ws1.Range("E1").Value2
)
$dataRow = @(
    $ws1.Range("A2").Value2,
    $ws1.Range("B2").Value2,
    $ws1.Range("C2").Value2,
    $ws1.Range("D2").Value2,
    $ws1.Range("E2").Value2
)

$ws1.Range("F1").Value = $headerRow[4]
$ws1.Range("E1").Value = $headerRow[3]
$ws1.Range("D1").Value = $headerRow[2]
$ws1.Range("C1").Value = $headerRow[1]
$ws1.Range("B1").Value = $headerRow[0]
$ws1.Range("A1").Value = "Test Scenario"

$ws1.Range("F2").Value = $dataRow[4]
$ws1.Range("E2").Value = $dataRow[3]
$ws1.Range("D2").Value = $dataRow[2]
$ws1.Range("C2").Value = $dataRow[1]
$ws1.Range("B2").Value = $dataRow[0]
$ws1.Range("A2").Value = "LoginPageTest"

# --- View state: select F3 on "Test Data" then make "Test Info" the active tab ---
$ws2.Range("F3").Select()
$ws1.Select()
